# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el texto de la conversión del día (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.21 = 8061.46 pesos`n✅ 8061.46 pesos = 2.2 = 941.42 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: actualizar las tasas en N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 451.99
$wsTasas.Range("O10").Value = 3643.7
$wsTasas.Range("N12").Value = 3665
$wsTasas.Range("O12").Value = 428
